# Insert a new data row before the current row 170 (shifts old rows 170-292 down to 171-293)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(170).Insert()

# Populate the newly inserted row 170 with the new record.
# Columns A,B,C,E,F,G,H,I,J,K,L,Q,R,T mirror the record that used to occupy row 170
# (now shifted to row 171); only D, M, N, O, P, S differ.
$ws.Range("A170").Value = 4
$ws.Range("B170").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C170").Value = "Los Lagos"
$ws.Range("D170").Value = 44879
$ws.Range("E170").Value = 10
$ws.Range("F170").Value = "Fruta"
$ws.Range("G170").Value = 100101
$ws.Range("H170").Value = "Berries"
$ws.Range("I170").Value = 100112025
$ws.Range("J170").Value = "Frutilla"
$ws.Range("K170").Value = "Sin especificar"
$ws.Range("L170").Value = "Primera"
$ws.Range("M170").Value = 600
$ws.Range("N170").Value = 10000
$ws.Range("O170").Value = 11000
$ws.Range("P170").Value = 10500
$ws.Range("Q170").Value = "`$/caja 7 kilos"
$ws.Range("R170").Value = "Región de La Araucanía"
$ws.Range("S170").Value = 1500
$ws.Range("T170").Value = 7
